$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update cell contents row by row to match the final data (typo fixes + new rows 18-24)

$ws.Range('A1').Value = 'Country'
$ws.Range('B1').Value = 'Species'
$ws.Range('C1').Value = 'Year'
$ws.Range('D1').Value = 'Source'
$ws.Range('E1').Value = 'GBADs Score'
$ws.Range('F1').Value = 'Reasoning'
$ws.Range('G1').Value = 'Sources'

$ws.Range('A2').Value = 'Chile'
$ws.Range('B2').Value = 'Cattle'
$ws.Range('C2').Value = 'All'
$ws.Range('D2').Value = 'WOAH'
$ws.Range('E2').Value = 1
$ws.Range('F2').Value = 'Over the 14 years that WOAH has reported data, there have only been two unique population numbers. It was 4757859 for two years, then 3719507 since then until 2018 before it abruptly ends. Although the population ranges fairly close to FAOSTAT between 2007 and 2012, given that the population numbers do not change at all, there is no way this is a reliable source for Chilean cattle population numbers'

$ws.Range('A3').Value = 'Chile'
$ws.Range('B3').Value = 'Cattle'
$ws.Range('C3').Value = 1961
$ws.Range('D3').Value = 'FAOSTAT'
$ws.Range('E3').Value = 5
$ws.Range('F3').Value = 'The population of this closely aligns with the population stated in the attached research paper which gives a high degrees of certainty that it is correct.'
$ws.Range('G3').Value = 'https://scholarworks.calstate.edu/downloads/cj82kb65d'

$ws.Range('A4').Value = 'Chile'
$ws.Range('B4').Value = 'Cattle'
$ws.Range('C4').Value = 2019
$ws.Range('D4').Value = 'FAOSTAT'
$ws.Range('E4').Value = 4
$ws.Range('F4').Value = ' The Chilean Agricultural Ministry (ODEPA)  put out a report in 2022 saying there was 3.108 million cattle in 2019, whereas FAOSTAT said it was 3.1 million. Since FAOSTAT and ODEPA have such similar numbers it is likely that FAOSTAT has good data. '
$ws.Range('G4').Value = 'https://bibliotecadigital.odepa.gob.cl/bitstream/handle/20.500.12650/69897/ficha_nacional_2023.pdf, https://www.odepa.gob.cl/wp-content/uploads/2019/09/panorama2019Final.pdf'

$ws.Range('A5').Value = 'Chile'
$ws.Range('B5').Value = 'Cattle'
$ws.Range('C5').Value = '1961-2000'
$ws.Range('D5').Value = 'FAOSTAT'
$ws.Range('E5').Value = 4
$ws.Range('F5').Value = "'All the data between these years (inclusive) all have 'official' flags attached to their data points. This means the points come from a government source in the country and the source has a high degree of confidence in the data points. "
$ws.Range('G5').Value = 'https://gbadske.org/dashboards/visualizer/'

$ws.Range('A6').Value = 'Chile'
$ws.Range('B6').Value = 'Cattle'
$ws.Range('C6').Value = 2008
$ws.Range('D6').Value = 'FAOSTAT'
$ws.Range('E6').Value = 5
$ws.Range('F6').Value = 'The Chilean Agricultural Ministry (ODEPA) put out a report in 2019 saying there was 3.79 million Cattle alive in 2008. FAOSTAT says there was 3.8 million cattle alive that year and has flagged it as official. This good evidence that FAOSTAT had good data for this year.'
$ws.Range('G6').Value = 'https://www.odepa.gob.cl/wp-content/uploads/2019/09/panorama2019Final.pdf'

$ws.Range('A7').Value = 'Chile'
$ws.Range('B7').Value = 'Pigs'
$ws.Range('C7').Value = 2019
$ws.Range('D7').Value = 'FAOSTAT'
$ws.Range('E7').Value = 5
$ws.Range('F7').Value = 'The Chilean Agricultural Ministry (ODEPA) put out a report in 2022 saying there was 2.767 million pigs alive in 2019. FAOSTAT says there was 2.57 million pigs alive that year and has flagged it as official. This good evidence that FAOSTAT had good data for this year.'
$ws.Range('G7').Value = 'https://bibliotecadigital.odepa.gob.cl/bitstream/handle/20.500.12650/69897/ficha_nacional_2023.pdf'

$ws.Range('A8').Value = 'Chile'
$ws.Range('B8').Value = 'Pigs'
$ws.Range('C8').Value = 2008
$ws.Range('D8').Value = 'FAOSTAT'
$ws.Range('E8').Value = 3
$ws.Range('F8').Value = 'ODEPA put out a report in 2019 saying there was almost 3.8 million Pigs alive in 2008. FAOSTAT reports there being 2.79 million alive at that time. That’s a difference of more than 25 percent. This is not good evidence for FAOSTAT as these is a pretty wide margin. '
$ws.Range('G8').Value = 'https://www.odepa.gob.cl/wp-content/uploads/2019/09/panorama2019Final.pdf'

$ws.Range('A9').Value = 'Chile '
$ws.Range('B9').Value = 'Pigs'
$ws.Range('C9').Value = 'all'
$ws.Range('D9').Value = 'FAOSTAT'
$ws.Range('E9').Value = 5
$ws.Range('F9').Value = 'The Chilean pork sector has focused more on exporting than trying to increase domestic consumption. Most of the pork exported goes to Asia and has allowed the pork sector to grow. '
$ws.Range('G9').Value = 'https://www.porkbusiness.com/news/industry/chile-demanding-pork-export-market-knows-what-it-needs'

$ws.Range('A10').Value = 'Chile'
$ws.Range('B10').Value = 'Chickens'
$ws.Range('C10').Value = '2013-2014'
$ws.Range('D10').Value = 'FAOSTAT'
$ws.Range('E10').Value = 0
$ws.Range('F10').Value = 'Impossible population jump. More than doubled the chicken population in a single year. This is most likely a similar thing to the USA chicken data were they switched from a single point in time census to a total yearly census. Does not align with other data collected from this time period. '
$ws.Range('G10').Value = 'https://www.statista.com/statistics/1002893/chile-chicken-meat-production-volume/'

$ws.Range('A11').Value = 'Chile'
$ws.Range('B11').Value = 'Chickens'
$ws.Range('C11').Value = '2016-2020'
$ws.Range('D11').Value = 'FAOSTAT'
$ws.Range('E11').Value = 4
$ws.Range('F11').Value = 'WATT Poultry international states that Chiles chicken population has grown from 240 million to 301 million chickens between 2016 and 2020. This conflicts with FAOSTAT as FAOSTAT says there is only 115 million chicken in Chile in 2020 and similar amounts in 2016. One possible explanation for this is the same issue that the United States Chicken population has where the data switch from being an instantaneous census to a  year long census. So WATT Poultry could be taking the a census of all chickens alive during this time. '
$ws.Range('G11').Value = 'https://www.wattagnet.com/blogs/blog/15534902/impact-of-chiles-new-government-on-animal-production#:~:text=The%20average%20Chilean%20consumes%2032.2,85%25%20of%20the%20region%27s%20average.'

$ws.Range('A12').Value = 'Chile'
$ws.Range('B12').Value = 'Chickens'
$ws.Range('C12').Value = '2002-2020'
$ws.Range('D12').Value = 'FAOSTAT'
$ws.Range('E12').Value = 5
$ws.Range('F12').Value = 'WATT Poultry International States that the chicken population in Chile has been growing because they have been exporting more chicken than they have been consuming. This is what is driving their chicken population to grow. This reasoning corroborates with the chicken population numbers FAOSTAT reports since it has been growing much faster over the last 2 decades. With most of their exporting''s going to Asia like their pork production.'
$ws.Range('G12').Value = 'https://www.wattagnet.com/blogs/blog/15534902/impact-of-chiles-new-government-on-animal-production#:~:text=The%20average%20Chilean%20consumes%2032.2,85%25%20of%20the%20region%27s%20average, https://www.statista.com/statistics/1002893/chile-chicken-meat-production-volume/, https://www.euromeatnews.com/Article-Chile-exported-637,000-tons-of-chicken,-turkey-and-pork-meat-in-2022/6200'

$ws.Range('A13').Value = 'Chile'
$ws.Range('B13').Value = 'Chickens'
$ws.Range('C13').Value = 2023
$ws.Range('D13').Value = 'All'
$ws.Range('E13').Value = 0
$ws.Range('F13').Value = 'Reuters claimed in an article that there was only 30 million chickens alive in Chile. This figure was from Carlos Orellana, head of livestock protection for Chile''s farming and livestock SAG agency who said it in a press conference in Santiago Chile. '
$ws.Range('G13').Value = 'https://www.reuters.com/business/healthcare-pharmaceuticals/chile-culls-40000-poultry-amid-industrial-bird-flu-outbreak-2023-03-15/'
$ws.Range('G13').Font.Name = 'Calibri'

$ws.Range('A14').Value = 'Chile'
$ws.Range('B14').Value = 'Sheep'
$ws.Range('C14').Value = 2018
$ws.Range('D14').Value = 'FAOSTAT'
$ws.Range('E14').Value = 3
$ws.Range('F14').Value = 'The attached paper cites the Instituto Nacional de Estadísticas (National Statistics Institute) which states there was 4 million sheep alive in 2018. FAOSTAT says was 2.1 million alive at that time'
$ws.Range('G14').Value = 'https://www.mdpi.com/2076-2615/8/10/165#:~:text=Introduction,in%20the%20form%20of%20tussocks).'

$ws.Range('A15').Value = 'Chile'
$ws.Range('B15').Value = 'Sheep'
$ws.Range('C15').Value = '2007-2018'
$ws.Range('D15').Value = 'WOAH'
$ws.Range('E15').Value = 0
$ws.Range('F15').Value = 'The population is the same for every year which is extremely unlikely'
$ws.Range('G15').Value = 'https://www.mdpi.com/2076-2615/8/10/165#:~:text=Introduction,in%20the%20form%20of%20tussocks).'

$ws.Range('A16').Value = 'Chile'
$ws.Range('B16').Value = 'Sheep'
$ws.Range('C16').Value = 'All'
$ws.Range('D16').Value = 'FAOSTAT'
$ws.Range('E16').Value = 3
$ws.Range('F16').Value = 'The population of sheep is kept mostly in the Pategonia region of Chile (Southern) in extremely large herds (thousands). This makes it hard to count the sheep because  '

$ws.Range('A17').Value = 'Chile'
$ws.Range('B17').Value = 'Sheep'
$ws.Range('C17').Value = 2004
$ws.Range('D17').Value = 'FAOSTAT'
$ws.Range('E17').Value = 3
$ws.Range('F17').Value = 'A paper by Carlos Alejandro Robles cites the Argentine wool federation who state in 2004 there was over 8 million sheep in Patagonia alone. This is double the number that FAOSTAT has for the same time. The Argentine Wool Association does not have an English translation on their website so it could not be verified at this time. '
$ws.Range('G17').Value = 'https://www.researchgate.net/figure/Numbers-of-sheep-and-wool-production-of-each-province-of-Patagonia-relative-to-national_tbl1_229637472'

$ws.Range('A18').Value = 'Chile'
$ws.Range('B18').Value = 'Cattle'
$ws.Range('C18').Value = '1971-1981'
$ws.Range('D18').Value = 'FAOSTAT'
$ws.Range('E18').Value = 4
$ws.Range('F18').Value = 'In the period between 1971-1981, inflation was extremely high in Chile reaching a peak of ~505% in 1974. This would have made it extremely hard for chileans to afford essential products because of the rampant inflation. However, Chile does not grow enough food to feed all its citizens and is reliant on imports from other other countries. Given the high inflation, this would have made it extremely expensive to buy imported foods and as a result kept their domestic food production high as that would have been more affordable that imported foods.'
$ws.Range('G18').Value = 'https://www.worlddata.info/america/chile/inflation-rates.php'

$ws.Range('A19').Value = 'Chile'
$ws.Range('B19').Value = 'All'
$ws.Range('C19').Value = '1973-1990'
$ws.Range('D19').Value = 'FAOSTAT'
$ws.Range('E19').Value = 5
$ws.Range('F19').Value = 'There was a armed resistance / civil war in Chile during this period where communist guerrillas attempted to establisht themselves and overthrow the government. However, they never got very large and didn''t have a massive impact on the country. This would have had minimal impact on the markets, population, and livestock.'
$ws.Range('G19').Value = 'https://www.theguardian.com/world/2023/sep/03/fifty-years-on-the-lasting-tragedy-of-chiles-coup'

$ws.Range('A20').Value = 'Chile'
$ws.Range('B20').Value = 'Cattle'
$ws.Range('C20').Value = '2010-2023'
$ws.Range('D20').Value = 'FAOSTAT'
$ws.Range('E20').Value = 5
$ws.Range('F20').Value = 'Starting in 2010, Chile went into a drought and it''s still occuring. It happened towards the south of the country which affect primarily the cattle and sheep farmers. The drought got progressivly worse until grass stop growing which is bad for farmers because this is the primary caloric input for these animals. As a result we would expect to see it impact the cattle population. After 2012 we see the cattle population dropping siginificantly. There could be other external factors impacting this decrease in population but it lines up with the drought occurance. '
$ws.Range('G20').Value = 'https://www.reuters.com/world/americas/the-grass-does-not-grow-chiles-far-south-worst-drought-50-years-2023-02-09/'

$ws.Range('A21').Value = 'Chile'
$ws.Range('B21').Value = 'Sheep'
$ws.Range('C21').Value = '2010-2023'
$ws.Range('D21').Value = 'FAOSTAT'
$ws.Range('E21').Value = 5
$ws.Range('F21').Value = 'Starting in 2010, Chile went into a drought and it''s still occuring. It happened towards the south of the country which affect primarily the cattle and sheep farmers. The drought got progressivly worse until grass stop growing which is bad for farmers because this is the primary caloric input for these animals. As a result we would expect to see it impact the sheep population. After 2012 we see the sheep population dropping siginificantly. There could be other external factors impacting this decrease in population but it lines up with the drought occurance. '
$ws.Range('G21').Value = 'https://www.reuters.com/world/americas/the-grass-does-not-grow-chiles-far-south-worst-drought-50-years-2023-02-09/'

$ws.Range('A22').Value = 'Chile'
$ws.Range('B22').Value = 'Cattle'
$ws.Range('C22').Value = '1967-1969'
$ws.Range('D22').Value = 'FAOSTAT'
$ws.Range('E22').Value = 2
$ws.Range('F22').Value = 'Between 1967 and 1969 there was a severe drought in Chile. Thousands of livestock died and farmers needed grants from the government to stay open. As a result we should see an impact on the livestock numbers over this period. We only see the population drop in 1971. A drop does make sense but the timing is off. We should see the population drop two years earlier. We see the population actually increase during the drought period even though sources say 100k+ cattle died during this time. '
$ws.Range('G22').Value = 'https://www.redalyc.org/journal/811/81172097005/html/'

$ws.Range('A23').Value = 'Chile'
$ws.Range('B23').Value = 'Sheep'
$ws.Range('C23').Value = '1967-1969'
$ws.Range('D23').Value = 'FAOSTAT'
$ws.Range('E23').Value = 5
$ws.Range('F23').Value = 'Between 1967 and 1969 there was a severe drought in Chile. Thousands of livestock died and farmers needed grants from the government to stay open. As a result we should see an impact on the livestock numbers over this period. After 1968 we start to see the sheep population start dropping. It drops for multiple years indicating a lasting impact from the drought. '
$ws.Range('G23').Value = 'https://www.redalyc.org/journal/811/81172097005/html/'

$ws.Range('A24').Value = 'Chile'

# Add hyperlink for the inflation-rates source cited in row 18
$ws.Hyperlinks.Add($ws.Range('G18'), 'https://www.worlddata.info/america/chile/inflation-rates.php')

# Update the view: selection moved to B24, zoomed in to 159%
$ws.Range('B24').Select()
$excel.ActiveWindow.Zoom = 159
